$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(4).Delete()
$ws.Range("L1:L28").Copy() | Out-Null
$ws.Range("M1:M28").PasteSpecial(-4122) | Out-Null
